# Updated cryptos list on Fri Apr 14 10:11:06 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.844.75"
$ws.Range("E2").Value = "  +2.02%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.114.58"
$ws.Range("E3").Value = "  +6.53%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.00"
$ws.Range("E5").Value = "  +3.16%  "

# Row 6 - USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5317"
$ws.Range("E7").Value = "  +4.13%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4382"
$ws.Range("E8").Value = "  +6.78%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09013"
$ws.Range("E9").Value = "  +5.05%  "

# Row 10 - OKB
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.26"
$ws.Range("E10").Value = "  +8.59%  "

# Row 11 - Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.181"
$ws.Range("E11").Value = "  +4.26%  "

# Row 12 - Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.11"
$ws.Range("E12").Value = "  +3.84%  "

# Row 13 - WrappedEther
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.115.83"
$ws.Range("E13").Value = "  +6.69%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.767"
$ws.Range("E14").Value = "  +4.43%  "

# Row 15 - Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.836"
$ws.Range("E15").Value = "  +6.02%  "

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.36"
$ws.Range("E16").Value = "  +3.66%  "

# Row 17 - BinanceUSD
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.10%  "

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001131"
$ws.Range("E18").Value = "  +1.76%  "

# Row 19 - TRON (price unchanged)
$ws.Range("E19").Value = "  +1.74%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.13"
$ws.Range("E20").Value = "  +1.56%  "

# Row 21 - Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.345"
$ws.Range("E22").Value = "  +4.06%  "

# Row 23 - WrappedBTC
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.892.04"
$ws.Range("E23").Value = "  +1.93%  "

# Row 24 - Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.42"
$ws.Range("E24").Value = "  +7.95%  "

# Row 25 - WrappedliquidstakedEther2.0
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.357.05"
$ws.Range("E25").Value = "  +6.53%  "

# Row 26 - Toncoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.269"
$ws.Range("E26").Value = "  +2.62%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.79"
$ws.Range("E27").Value = "  +1.55%  "

# Row 28 - LidoDAOToken (price unchanged)
$ws.Range("E28").Value = "  +8.48%  "

# Row 29 - Monero
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.33"
$ws.Range("E29").Value = "  +0.24%  "

# Row 30 - BitcoinCash
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.58"
$ws.Range("E30").Value = "  +2.04%  "

# Row 31 - ImmutableX
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.187"
$ws.Range("E31").Value = "  +4.47%  "

# Row 32 - Stellar (price unchanged)
$ws.Range("E32").Value = "  +2.14%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.244"
$ws.Range("E33").Value = "  +3.45%  "

# Row 34 - HuobiToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.014"
$ws.Range("E34").Value = "  +5.42%  "

# Row 35 - ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.564"
$ws.Range("E35").Value = "  +18.49%  "

# Row 36 - VeChain
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02623"
$ws.Range("E36").Value = "  +5.60%  "

# Row 37 - Aptos
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.91"
$ws.Range("E37").Value = "  +8.96%  "

# Row 38 - InternetComputer(DFINITY)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.542"
$ws.Range("E38").Value = "  +2.78%  "

# Row 39 - Hedera
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06763"
$ws.Range("E39").Value = "  +4.13%  "

# Row 40 - FraxShare (price unchanged)
$ws.Range("E40").Value = "  +5.88%  "

# Row 41 - Algorand
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2288"
$ws.Range("E41").Value = "  +5.06%  "

# Row 42 - TheSandbox
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6859"
$ws.Range("E42").Value = "  +4.34%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.256"
$ws.Range("E43").Value = "  +2.91%  "

# Row 44 - was EnergySwap, now becomes Decentraland (row 44/45 content swap)
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6478"
$ws.Range("E44").Value = "  +5.84%  "

# Row 45 - was Decentraland, now becomes EnergySwap
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.14"
$ws.Range("E45").Value = "  +3.97%  "

# Row 46 - Frax
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9992"
$ws.Range("E46").Value = "  -0.02%  "

# Row 47 - NEARProtocol
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.232"
$ws.Range("E47").Value = "  +1.93%  "

# Row 48 - PancakeSwap
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.665"
$ws.Range("E48").Value = "  +0.36%  "

# Row 49 - EOS (price unchanged)
$ws.Range("E49").Value = "  +4.48%  "

# Row 50 - Aave
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.99"
$ws.Range("E50").Value = "  +4.15%  "

# Row 51 - Quant
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.53"
$ws.Range("E51").Value = "  -1.84%  "
